$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Twitter account name / contributor string in cell C6
$ws.Range("C6").Value = "Joan Martinez (Twitter @jjenifer457)"
